# Update "想去人数" (F column) values on the "展览" and "全部类型" worksheets.
# Both sheets share the same first 5 rows of data, so the same F2:F5 updates
# apply to each.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1447
    $ws.Range("F3").Value = 3054
    $ws.Range("F4").Value = 37
    $ws.Range("F5").Value = 538
}
